# Update TPM-derived metric values in LR-pairs sheet (Il1b-Il1r1) per new TPM data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.420398"
$ws.Range("H2").Value = [double]"1.261194"
$ws.Range("I2").Value = [double]"0.001794717364332138"
$ws.Range("J2").Value = [double]"0.001794717364332138"
$ws.Range("M2").Value = [double]"12.673913"
$ws.Range("N2").Value = [double]"38.021739"
$ws.Range("O2").Value = [double]"0.234043494199914"
$ws.Range("P2").Value = [double]"0.234043494199914"
$ws.Range("Q2").Value = [double]"5.328087677374"
$ws.Range("R2").Value = [double]"47.952789096366"
$ws.Range("S2").Value = [double]"0.0004200419230495537"
$ws.Range("T2").Value = [double]"0.0004200419230495536"
$ws.Range("G3").Value = [double]"0.420398"
$ws.Range("H3").Value = [double]"1.261194"
$ws.Range("I3").Value = [double]"0.001794717364332138"
$ws.Range("J3").Value = [double]"0.001794717364332138"
$ws.Range("O3").Value = [double]"0.6388539132363013"
$ws.Range("P3").Value = [double]"0.6388539132363011"
$ws.Range("Q3").Value = [double]"14.543748265222"
$ws.Range("R3").Value = [double]"130.893734386998"
$ws.Range("S3").Value = [double]"0.001146562211356727"
$ws.Range("T3").Value = [double]"0.001146562211356727"
$ws.Range("G4").Value = [double]"0.420398"
$ws.Range("H4").Value = [double]"1.261194"
$ws.Range("I4").Value = [double]"0.001794717364332138"
$ws.Range("J4").Value = [double]"0.001794717364332138"
$ws.Range("M4").Value = [double]"6.728406666666667"
$ws.Range("N4").Value = [double]"20.18522"
$ws.Range("O4").Value = [double]"0.1242504825987572"
$ws.Range("P4").Value = [double]"0.1242504825987572"
$ws.Range("Q4").Value = [double]"2.828608705853334"
$ws.Range("R4").Value = [double]"25.45747835268"
$ws.Range("S4").Value = [double]"0.0002229944986466377"
$ws.Range("T4").Value = [double]"0.0002229944986466377"
$ws.Range("G5").Value = [double]"0.420398"
$ws.Range("H5").Value = [double]"1.261194"
$ws.Range("I5").Value = [double]"0.001794717364332138"
$ws.Range("J5").Value = [double]"0.001794717364332138"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.1544473333333333"
$ws.Range("N5").Value = [double]"0.463342"
$ws.Range("O5").Value = [double]"0.002852109965027549"
$ws.Range("P5").Value = [double]"0.002852109965027548"
$ws.Range("Q5").Value = [double]"0.06492935003866668"
$ws.Range("R5").Value = [double]"0.584364150348"
$ws.Range("S5").Value = [double]"5.118731279219667E-06"
$ws.Range("T5").Value = [double]"5.118731279219666E-06"
$ws.Range("G6").Value = [double]"0.04680500000000001"
$ws.Range("H6").Value = [double]"0.140415"
$ws.Range("I6").Value = [double]"0.0001998148093891163"
$ws.Range("J6").Value = [double]"0.0001998148093891163"
$ws.Range("M6").Value = [double]"12.673913"
$ws.Range("N6").Value = [double]"38.021739"
$ws.Range("O6").Value = [double]"0.234043494199914"
$ws.Range("P6").Value = [double]"0.234043494199914"
$ws.Range("Q6").Value = [double]"0.593202497965"
$ws.Range("R6").Value = [double]"5.338822481685"
$ws.Range("S6").Value = [double]"4.676535618231857E-05"
$ws.Range("T6").Value = [double]"4.676535618231855E-05"
$ws.Range("G7").Value = [double]"0.04680500000000001"
$ws.Range("H7").Value = [double]"0.140415"
$ws.Range("I7").Value = [double]"0.0001998148093891163"
$ws.Range("J7").Value = [double]"0.0001998148093891163"
$ws.Range("O7").Value = [double]"0.6388539132363013"
$ws.Range("P7").Value = [double]"0.6388539132363011"
$ws.Range("Q7").Value = [double]"1.619227821145"
$ws.Range("R7").Value = [double]"14.573050390305"
$ws.Range("S7").Value = [double]"0.0001276524729008026"
$ws.Range("T7").Value = [double]"0.0001276524729008025"
$ws.Range("G8").Value = [double]"0.04680500000000001"
$ws.Range("H8").Value = [double]"0.140415"
$ws.Range("I8").Value = [double]"0.0001998148093891163"
$ws.Range("J8").Value = [double]"0.0001998148093891163"
$ws.Range("M8").Value = [double]"6.728406666666667"
$ws.Range("N8").Value = [double]"20.18522"
$ws.Range("O8").Value = [double]"0.1242504825987572"
$ws.Range("P8").Value = [double]"0.1242504825987572"
$ws.Range("Q8").Value = [double]"0.3149230740333334"
$ws.Range("R8").Value = [double]"2.8343076663"
$ws.Range("S8").Value = [double]"2.482708649697639E-05"
$ws.Range("T8").Value = [double]"2.482708649697638E-05"
$ws.Range("G9").Value = [double]"0.04680500000000001"
$ws.Range("H9").Value = [double]"0.140415"
$ws.Range("I9").Value = [double]"0.0001998148093891163"
$ws.Range("J9").Value = [double]"0.0001998148093891163"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.1544473333333333"
$ws.Range("N9").Value = [double]"0.463342"
$ws.Range("O9").Value = [double]"0.002852109965027549"
$ws.Range("P9").Value = [double]"0.002852109965027548"
$ws.Range("Q9").Value = [double]"0.007228907436666667"
$ws.Range("R9").Value = [double]"0.06506016693"
$ws.Range("S9").Value = [double]"5.698938090187787E-07"
$ws.Range("T9").Value = [double]"5.698938090187785E-07"
$ws.Range("G10").Value = [double]"233.774694"
$ws.Range("H10").Value = [double]"701.324082"
$ws.Range("I10").Value = [double]"0.9980054678262787"
$ws.Range("J10").Value = [double]"0.9980054678262787"
$ws.Range("M10").Value = [double]"12.673913"
$ws.Range("N10").Value = [double]"38.021739"
$ws.Range("O10").Value = [double]"0.234043494199914"
$ws.Range("P10").Value = [double]"0.234043494199914"
$ws.Range("Q10").Value = [double]"2962.840133357622"
$ws.Range("R10").Value = [double]"26665.5612002186"
$ws.Range("S10").Value = [double]"0.2335766869206822"
$ws.Range("T10").Value = [double]"0.2335766869206821"
$ws.Range("G11").Value = [double]"233.774694"
$ws.Range("H11").Value = [double]"701.324082"
$ws.Range("I11").Value = [double]"0.9980054678262787"
$ws.Range("J11").Value = [double]"0.9980054678262787"
$ws.Range("O11").Value = [double]"0.6388539132363013"
$ws.Range("P11").Value = [double]"0.6388539132363011"
$ws.Range("Q11").Value = [double]"8087.479722347165"
$ws.Range("R11").Value = [double]"72787.31750112449"
$ws.Range("S11").Value = [double]"0.6375796985520438"
$ws.Range("T11").Value = [double]"0.6375796985520437"
$ws.Range("G12").Value = [double]"233.774694"
$ws.Range("H12").Value = [double]"701.324082"
$ws.Range("I12").Value = [double]"0.9980054678262787"
$ws.Range("J12").Value = [double]"0.9980054678262787"
$ws.Range("M12").Value = [double]"6.728406666666667"
$ws.Range("N12").Value = [double]"20.18522"
$ws.Range("O12").Value = [double]"0.1242504825987572"
$ws.Range("P12").Value = [double]"0.1242504825987572"
$ws.Range("Q12").Value = [double]"1572.93120960756"
$ws.Range("R12").Value = [double]"14156.38088646804"
$ws.Range("S12").Value = [double]"0.1240026610136136"
$ws.Range("T12").Value = [double]"0.1240026610136136"
$ws.Range("G13").Value = [double]"233.774694"
$ws.Range("H13").Value = [double]"701.324082"
$ws.Range("I13").Value = [double]"0.9980054678262787"
$ws.Range("J13").Value = [double]"0.9980054678262787"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"0.6666666666666666"
$ws.Range("M13").Value = [double]"0.1544473333333333"
$ws.Range("N13").Value = [double]"0.463342"
$ws.Range("O13").Value = [double]"0.002852109965027549"
$ws.Range("P13").Value = [double]"0.002852109965027548"
$ws.Range("Q13").Value = [double]"36.105878089116"
$ws.Range("R13").Value = [double]"324.952902802044"
$ws.Range("S13").Value = [double]"0.00284642133993931"
$ws.Range("T13").Value = [double]"0.00284642133993931"
